$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.114.08'
$ws.Range('E2').Value = '  -3.26%  '
$ws.Range('D3').Value = '2.521.46'
$ws.Range('E3').Value = '  -4.55%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '''577.66'
$ws.Range('E5').Value = '  -1.39%  '
$ws.Range('D6').Value = '''168.15'
$ws.Range('E6').Value = '  -4.20%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '2.520.56'
$ws.Range('E9').Value = '  -4.49%  '
$ws.Range('D10').Value = '''0.162'
$ws.Range('E10').Value = '  -5.70%  '
$ws.Range('E11').Value = '  -1.66%  '
$ws.Range('D12').Value = '''0.347'
$ws.Range('E12').Value = '  -3.50%  '
$ws.Range('D13').Value = '''4.90'
$ws.Range('E13').Value = '  -0.79%  '
$ws.Range('D14').Value = '2.981.13'
$ws.Range('E14').Value = '  -4.71%  '
$ws.Range('D15').Value = '69.979.95'
$ws.Range('E15').Value = '  -3.21%  '
$ws.Range('D16').Value = '''0.0000176'
$ws.Range('E16').Value = '  -5.68%  '
$ws.Range('D17').Value = '''25.13'
$ws.Range('E17').Value = '  -2.55%  '
$ws.Range('D18').Value = '2.532.41'
$ws.Range('E18').Value = '  -4.03%  '
$ws.Range('D19').Value = '''7.78'
$ws.Range('E19').Value = '  -0.93%  '
$ws.Range('D20').Value = '''11.34'
$ws.Range('E20').Value = '  -6.03%  '
$ws.Range('D21').Value = '''351.61'
$ws.Range('E21').Value = '  -6.43%  '
$ws.Range('E22').Value = '  -4.00%  '
$ws.Range('E23').Value = '  -3.76%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = '''69.16'
$ws.Range('E25').Value = '  -3.31%  '
$ws.Range('D26').Value = '''4.02'
$ws.Range('E26').Value = '  -5.33%  '
$ws.Range('D27').Value = '''9.03'
$ws.Range('E27').Value = '  -5.06%  '
$ws.Range('E28').Value = '  -4.72%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').Value = '0.0₃0911'
$ws.Range('E30').Value = '  -4.29%  '
$ws.Range('D31').Value = '''7.91'
$ws.Range('E31').Value = '  -0.87%  '
$ws.Range('E32').Value = '  -2.79%  '
$ws.Range('D33').Value = '''466.17'
$ws.Range('E33').Value = '  -5.29%  '
$ws.Range('E34').Value = '  -2.09%  '
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').Value = '''0.121'
$ws.Range('E36').Value = '  +4.43%  '
$ws.Range('D37').Value = '''153.01'
$ws.Range('E37').Value = '  -5.37%  '
$ws.Range('E38').Value = '  +0.76%  '
$ws.Range('E39').Value = '  -3.67%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').Value = '''4.79'
$ws.Range('E41').Value = '  -2.14%  '
$ws.Range('E42').Value = '  -1.69%  '
$ws.Range('E43').Value = '  -6.78%  '
$ws.Range('E44').Value = '  -13.99%  '
$ws.Range('D45').Value = '''2.33'
$ws.Range('E45').Value = '  -9.82%  '
$ws.Range('D46').Value = '''38.21'
$ws.Range('E46').Value = '  -2.14%  '
$ws.Range('D47').Value = '''143.59'
$ws.Range('E47').Value = '  -4.65%  '
$ws.Range('E48').Value = '  -2.02%  '
$ws.Range('E49').Value = '  -3.36%  '
$ws.Range('E50').Value = '  -4.63%  '
$ws.Range('D51').Value = '''0.0735'
$ws.Range('E51').Value = '  -1.30%  '
